# Atualização automática dos dados do dashboard (aba "Entrada")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entrada")

# Linha 3 - FERRAMENTAS/ MATRIZARIA
$ws.Range("B3").Value = "R$ 541.715,15"
$ws.Range("D3").Value = "R$ 541.715,15"
$ws.Range("F3").Value = "79,66 %"

# Linha 4 - MATERIA PRIMA
$ws.Range("B4").Value = "R$ 371.686,93"
$ws.Range("D4").Value = "R$ 371.686,93"
$ws.Range("F4").Value = "37.168.693,00 %"

# Linha 5 - agora CUSTO DESENVOLVIMENTO (era REFUGO REAL (PROCESSO))
$ws.Range("A5").Value = "CUSTO DESENVOLVIMENTO"
$ws.Range("B5").Value = "R$ 370.293,54"
$ws.Range("D5").Value = "R$ 370.293,54"
$ws.Range("E5").Value = "R$ 370.293,54"

# Linha 6 - agora REFUGO REAL (PROCESSO) (era CUSTO DESENVOLVIMENTO)
$ws.Range("A6").Value = "REFUGO REAL (PROCESSO)"
$ws.Range("B6").Value = "R$ 357.461,99"
$ws.Range("D6").Value = "R$ 357.461,99"
$ws.Range("E6").Value = "R$ 357.461,99"

# Linha 7 - FRETES
$ws.Range("B7").Value = "R$ 258.961,64"
$ws.Range("D7").Value = "R$ 258.961,64"
$ws.Range("F7").Value = "68,87 %"

# Linha 8 - agora REFUGO MP+CP* (era MANUTENCAO)
$ws.Range("A8").Value = "REFUGO MP+CP*"
$ws.Range("B8").Value = "R$ 212.131,03"
$ws.Range("C8").Value = "R$ 0,00"
$ws.Range("D8").Value = "R$ 212.131,03"
$ws.Range("E8").Value = "R$ 280.000,00"
$ws.Range("F8").Value = "75,76 %"

# Linha 9 - agora MANUTENCAO (era REFUGO MP+CP*)
$ws.Range("A9").Value = "MANUTENCAO"
$ws.Range("B9").Value = "R$ 209.640,19"
$ws.Range("C9").Value = "R$ 192.005,09"
$ws.Range("D9").Value = "R$ 401.645,28"
$ws.Range("E9").Value = "R$ 480.000,00"
$ws.Range("F9").Value = "83,68 %"

# Linha 10 - DESP. INDUSTRIAL
$ws.Range("B10").Value = "R$ 168.751,13"
$ws.Range("C10").Value = "R$ 121.898,00"
$ws.Range("D10").Value = "R$ 290.649,13"
$ws.Range("F10").Value = "64,59 %"

# Linha 11 - OLEOS E LUBRIFICANTES
$ws.Range("B11").Value = "R$ 110.014,09"
$ws.Range("D11").Value = "R$ 218.781,75"
$ws.Range("F11").Value = "78,14 %"

# Linha 12 - EMBALAGENS
$ws.Range("C12").Value = "R$ 74.942,31"
$ws.Range("D12").Value = "R$ 133.506,91"
$ws.Range("F12").Value = "83,44 %"

# Linha 13 - SERVICOS DE QUALIDADE
$ws.Range("B13").Value = "R$ 49.915,98"
$ws.Range("D13").Value = "R$ 49.915,98"
$ws.Range("F13").Value = "110,92 %"

# Linha 14 - FERRAMENTARIA/MAN FR
$ws.Range("B14").Value = "R$ 11.273,87"
$ws.Range("D14").Value = "R$ 11.273,87"
$ws.Range("F14").Value = "32,21 %"

# Linha 15 - CUSTO COM DESENVOLVIMENTO
$ws.Range("B15").Value = "R$ 9.192,99"
$ws.Range("D15").Value = "R$ 9.192,99"
$ws.Range("E15").Value = "R$ 9.192,99"

# Linha 18 - Total Geral
$ws.Range("B18").Value = "R$ 3.294.708,58"
$ws.Range("C18").Value = "R$ 502.145,83"
$ws.Range("D18").Value = "R$ 3.796.854,41"
$ws.Range("E18").Value = "R$ 4.584.075,63"
$ws.Range("F18").Value = "82,83 %"
